$wb = $excel.ActiveWorkbook

$wsValid = $wb.Worksheets.Item("validLoginData")
$wsInvalid = $wb.Worksheets.Item("invalidLoginData")

# --- validLoginData sheet (sql "valid" data -> now points at the admin DB creds) ---
$wsValid.Range("A2:B2").Style = "Normal"
$wsValid.Range("A3").Style = "Normal"
$wsValid.Range("A3").Value = "adminadmin"

$wsValid.Range("B3").Value = "Emre@Furkan28"
$hyperlink = $wsValid.Hyperlinks.Add($wsValid.Range("B3"), "mailto:Emre@Furkan28", "", "", "Emre@Furkan28")
$wsValid.Range("B3").Font.Name = "Arial"
$wsValid.Range("B3").Font.Size = 10
$wsValid.Range("B3").Font.Underline = 0
$wsValid.Range("B3").Font.Color = 16711680
$wsValid.Range("B3").Font.Charset = 1

# --- invalidLoginData sheet: values stay the same, just an explicit restyle ---
$wsInvalid.Range("A1:B3").Style = "Normal"

# --- switch the active tab back to validLoginData ---
$wsValid.Activate()
$null = $wsValid.Range("C6").Select()
